$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(39.06, 0.04000000000000001, 3.896192789077759, 1.09375),
    @(39.78, 0.16, 10.04770517349243, 7.53125),
    @(39.78, 0.16, 23.37266087532043, 14.328125),
    @(39.78, 0.16, 23.50989937782288, 14.453125),
    @(54.5, 0.16, 30.71380376815796, 18.90625),
    @(54.5, 0.16, 39.19084405899048, 17.96875),
    @(39.06, 0.04000000000000001, 6.48872709274292, 3.015625),
    @(3.6, 0.04000000000000001, 7.816614866256714, 3.609375),
    @(3.6, 0.04000000000000001, 4.427245378494263, 2.328125),
    @(59.83, 0.04000000000000001, 8.842975854873657, 3.484375),
    @(39.06, 0.04000000000000001, 103.8992028236389, 4.765625),
    @(39.06, 0.04000000000000001, 7.417251348495483, 2.984375),
    @(39.78, 0.16, 22.68746423721313, 13.84375),
    @(39.78, 0.16, 43.00482630729675, 18.4375),
    @(39.06, 0.04000000000000001, 6.663230180740356, 3.4375),
    @(39.06, 0.04000000000000001, 7.120179891586304, 3.390625),
    @(39.78, 0.16, 22.58138370513916, 14.84375),
    @(39.06, 0.04000000000000001, 6.986309766769409, 3.09375),
    @(39.06, 0.04000000000000001, 6.389602184295654, 3.15625),
    @(39.06, 0.04000000000000001, 6.762027502059937, 3.125),
    @(39.78, 0.16, 21.98983311653137, 15.828125),
    @(39.78, 0.16, 24.83666658401489, 9.46875),
    @(39.78, 0.16, 23.77997732162476, 11)
)

$startRow = 192
for ($i = 0; $i -lt $values.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $values[$i][0]
    $ws.Cells.Item($r, 2).Value = $values[$i][1]
    $ws.Cells.Item($r, 3).Value = $values[$i][2]
    $ws.Cells.Item($r, 4).Value = $values[$i][3]
}
